$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 171, pushing the existing data (rows 171-380)
# down by one, growing the used range to A1:R381.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44897
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112039
$ws.Range("G171").Value = "Ciboulette"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 160
$ws.Range("K171").Value = 1500
$ws.Range("L171").Value = 1500
$ws.Range("M171").Value = 1500
$ws.Range("N171").Value = "$/docena de atados"
$ws.Range("O171").Value = "Provincia de Quillota"
$ws.Range("P171").Value = 500
$ws.Range("Q171").Value = 3
$ws.Range("R171").Value = "Hortaliza"
